# Apply the "cryptos list" refresh: per-cell Price (D) / Volume(1h) (E) updates,
# plus a Filecoin/ImmutableX row swap (rows 31-32, columns B/C/D/E).
#
# Many Price values look like plain numbers ("1.00", "4.52", "61.35", ...).
# A naive `.Value = "1.00"` lets Excel auto-coerce that to the number 1, which
# would store it as a numeric cell instead of text. We force text by prefixing
# the literal with a leading apostrophe (the classic "treat as text" trick),
# then reset `.Style` back to "Normal" so no stray cell-style index is left
# behind on the cell (keeps the XML diff minimal / matches the source).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($cellRef).Style = "Normal"
}

# Row 2
Set-TextCell 'D2' '37.084.91'
Set-TextCell 'E2' '  +0.93%  '

# Row 3
Set-TextCell 'D3' '2.046.91'
Set-TextCell 'E3' '  -3.54%  '

# Row 4
Set-TextCell 'E4' '  -0.07%  '

# Row 5
Set-TextCell 'D5' '247.68'
Set-TextCell 'E5' '  -3.01%  '

# Row 6
Set-TextCell 'E6' '  -2.36%  '

# Row 7
Set-TextCell 'D7' '1.00'

# Row 8
Set-TextCell 'D8' '54.70'
Set-TextCell 'E8' '  +15.08%  '

# Row 9
Set-TextCell 'D9' '61.35'
Set-TextCell 'E9' '  +0.36%  '

# Row 10
Set-TextCell 'D10' '0.375'
Set-TextCell 'E10' '  +0.20%  '

# Row 11
Set-TextCell 'D11' '0.0760'
Set-TextCell 'E11' '  +2.63%  '

# Row 12
Set-TextCell 'E12' '  +5.06%  '

# Row 13
Set-TextCell 'D13' '14.97'
Set-TextCell 'E13' '  +3.60%  '

# Row 14
Set-TextCell 'D14' '2.339.49'
Set-TextCell 'E14' '  -3.96%  '

# Row 15
Set-TextCell 'D15' '0.812'
Set-TextCell 'E15' '  -3.95%  '

# Row 16
Set-TextCell 'D16' '5.17'
Set-TextCell 'E16' '  +0.70%  '

# Row 17
Set-TextCell 'D17' '2.040.07'
Set-TextCell 'E17' '  -3.91%  '

# Row 18
Set-TextCell 'D18' '37.017.48'
Set-TextCell 'E18' '  +0.80%  '

# Row 19
Set-TextCell 'D19' '71.74'
Set-TextCell 'E19' '  -2.76%  '

# Row 20
Set-TextCell 'D20' '0.0₃0895'
Set-TextCell 'E20' '  +6.35%  '

# Row 21
Set-TextCell 'D21' '14.11'
Set-TextCell 'E21' '  +5.87%  '

# Row 22
Set-TextCell 'D22' '235.85'
Set-TextCell 'E22' '  -2.37%  '

# Row 23
Set-TextCell 'D23' '5.23'
Set-TextCell 'E23' '  +0.12%  '

# Row 24
Set-TextCell 'E24' '  +0.09%  '

# Row 25
Set-TextCell 'D25' '2.40'
Set-TextCell 'E25' '  -3.31%  '

# Row 26
Set-TextCell 'D26' '169.18'
Set-TextCell 'E26' '  -1.56%  '

# Row 27
Set-TextCell 'D27' '8.97'
Set-TextCell 'E27' '  -2.27%  '

# Row 28
Set-TextCell 'D28' '19.95'
Set-TextCell 'E28' '  -7.51%  '

# Row 29
Set-TextCell 'E29' '  -2.37%  '

# Row 30
Set-TextCell 'E30' '  -1.15%  '

# Row 31
Set-TextCell 'B31' 'ImmutableX'
Set-TextCell 'C31' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell 'D31' '1.06'
Set-TextCell 'E31' '  +12.80%  '

# Row 32
Set-TextCell 'B32' 'Filecoin'
Set-TextCell 'C32' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell 'D32' '4.52'
Set-TextCell 'E32' '  +0.09%  '

# Row 33
Set-TextCell 'D33' '0.0617'
Set-TextCell 'E33' '  +2.67%  '

# Row 34
Set-TextCell 'D34' '4.32'
Set-TextCell 'E34' '  +3.28%  '

# Row 35
Set-TextCell 'E35' '  -0.08%  '

# Row 36
Set-TextCell 'D36' '0.0869'
Set-TextCell 'E36' '  -9.52%  '

# Row 37
Set-TextCell 'D37' '2.25'
Set-TextCell 'E37' '  -4.31%  '

# Row 38
Set-TextCell 'D38' '1.76'
Set-TextCell 'E38' '  -6.72%  '

# Row 39
Set-TextCell 'D39' '1.33'
Set-TextCell 'E39' '  -2.03%  '

# Row 40
Set-TextCell 'D40' '0.104'
Set-TextCell 'E40' '  +22.89%  '

# Row 41
Set-TextCell 'D41' '18.25'
Set-TextCell 'E41' '  +12.03%  '

# Row 42
Set-TextCell 'D42' '15.62'
Set-TextCell 'E42' '  -45.63%  '

# Row 43
Set-TextCell 'D43' '0.0222'
Set-TextCell 'E43' '  -1.73%  '

# Row 44
Set-TextCell 'E44' '  -5.22%  '

# Row 45
Set-TextCell 'D45' '95.18'
Set-TextCell 'E45' '  -4.06%  '

# Row 46
Set-TextCell 'D46' '2.78'
Set-TextCell 'E46' '  -0.99%  '

# Row 47
Set-TextCell 'D47' '4.04'
Set-TextCell 'E47' '  +40.24%  '

# Row 48
Set-TextCell 'D48' '1.288.54'
Set-TextCell 'E48' '  -5.06%  '

# Row 49
Set-TextCell 'E49' '  +2.06%  '

# Row 50
Set-TextCell 'E50' '  +1.51%  '

# Row 51
Set-TextCell 'D51' '6.71'
Set-TextCell 'E51' '  -5.49%  '
